$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (shifts the old rows 6-8 down to 7-9) and
# populate it with the Tokyo Nova: The 2nd Edition: Handbook entry.
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = 1996
$ws.Range("E6").Value = "tokyo_nova_handbook.jpg"
$ws.Range("B6").Value = "トーキョーN◎VA The 2nd Edition ハンドブック"
$ws.Range("C6").Value = "Tokyo Nova: The 2nd Edition: Handbook"
$ws.Range("D6").Value = "Aspect"
$ws.Range("F6").Value = "supplement"

# Update the selection to match the saved workbook state.
$ws.Range("C7").Select()
